$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Row 4
$ws.Range("K4").Value = 419.13

# Row 14
$ws.Range("K14").Value = 855.36
$ws.Range("L14").Value = 2266.66

# Row 16
$ws.Range("L16").Value = 0

# Row 18
$ws.Range("L18").Value = 6725.74

# Row 19
$ws.Range("D19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0

# Row 21
$ws.Range("C21").Value = 513.22
$ws.Range("M21").Value = 423.14
$ws.Range("N21").Value = 1058.37

# Row 29 (summary counts "x de 27")
$ws.Range("C29").Value = "1 de 27"
$ws.Range("D29").Value = "0 de 27"
$ws.Range("K29").Value = "2 de 27"
$ws.Range("M29").Value = "1 de 27"
$ws.Range("N29").Value = "1 de 27"
